$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.776.07'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '2.905.59'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '528.97'
$ws.Range('E5').Value = '  -2.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.08'
$ws.Range('E6').Value = '  -4.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  -2.28%  '
$ws.Range('D9').Value = '2.915.36'
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('E10').Value = '  -3.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.05'
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('E12').Value = '  -1.56%  '
$ws.Range('D13').Value = '3.415.24'
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.127'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '60.714.45'
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.87'
$ws.Range('E16').Value = '  -3.30%  '
$ws.Range('D17').Value = '2.903.73'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('E18').Value = '  -3.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.04'
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.74'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '362.83'
$ws.Range('E21').Value = '  -4.83%  '
$ws.Range('E22').Value = '  -0.46%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.69'
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.83'
$ws.Range('E25').Value = '  -0.69%  '
$ws.Range('E26').Value = '  -2.68%  '
$ws.Range('E27').Value = '  -3.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.86'
$ws.Range('E29').Value = '  -5.69%  '
$ws.Range('E30').Value = '  -7.79%  '
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('E33').Value = '  -2.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '151.98'
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('E35').Value = '  -5.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.59'
$ws.Range('E36').Value = '  -5.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -5.41%  '
$ws.Range('E38').Value = '  -4.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.88'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('E40').Value = '  -3.98%  '
$ws.Range('E41').Value = '  -5.07%  '
$ws.Range('D42').Value = '2.300.38'
$ws.Range('E42').Value = '  -4.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.651'
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0587'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('E45').Value = '  -7.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.00'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('E48').Value = '  -3.05%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0926'
$ws.Range('E49').Value = '  -3.21%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.31'
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '251.67'
$ws.Range('E51').Value = '  -6.10%  '
